# Add a new "price.value" column to the Item sheet (3rd sheet), between
# price.offered_value (I) and category_ids[0] (old J, now K), and make the
# Item sheet the active/selected sheet (mirrors "Added phonenumber and
# fixed template" commit).

$wb = $excel.ActiveWorkbook

$item = $wb.Worksheets.Item(3)

# Insert a new column before column J (10) on the Item sheet, shifting the
# existing price.value-and-beyond columns one to the right.
$item.Columns.Item(10).Insert()

# Populate the new column's header + data.
$item.Range("J1").Value = "price.value"
$item.Range("J2").Value = 80

# Make "Item" the active sheet/tab (workbook activeTab -> 2, Item tabSelected).
$item.Activate()
$item.Range("J1").Select()

Write-Host "edit applied"
